# Apply targeted cell fixes (Q58, Q64, Q72, O1140, R1142, R1143)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(58, 17).Value = 0
$ws.Cells.Item(64, 17).Value = 0
$ws.Cells.Item(72, 17).Value = 0
$ws.Cells.Item(1140, 15).Value = 2
$ws.Cells.Item(1142, 18).Value = 0
$ws.Cells.Item(1143, 18).Value = 0

# Append new weekly rows 1144-1165 (columns A-Q); column R left blank as in source data
$newRows = @(
    @(45474,114,131.5,111.8000030517578,130.3600006103516,130.1988677978516,251388266,2024,7,1,0,0,0,27,0,0,0),
    @(45481,131,135.3999938964844,120,120.6500015258789,120.5008773803711,186437818,2024,7,8,0,0,0,28,0,0,0),
    @(45488,121,124.1800003051758,113.8499984741211,115.5,115.3572387695312,72649724,2024,7,15,0,0,0,29,0,0,1),
    @(45495,113.370002746582,125.0999984741211,105.4000015258789,123.25,123.0976638793945,156855726,2024,7,22,0,0,0,30,2,0,0),
    @(45502,123.8099975585938,140.5,123.5,132.7700042724609,132.6058959960938,310188154,2024,7,29,0,0,0,31,0,0,0),
    @(45509,127.0999984741211,142.5,122.2699966430664,132.5099945068359,132.3462066650391,322394764,2024,8,5,0,0,0,32,0,0,0),
    @(45516,133.1000061035156,146.1999969482422,130.5,131.7100067138672,131.5472106933594,245135052,2024,8,12,0,0,0,33,0,0,0),
    @(45523,132.75,151.2299957275391,131.7200012207031,149.9799957275391,149.7946166992188,286202200,2024,8,19,0,0,0,34,0,0,0),
    @(45530,150.5,152,141.75,149.6199951171875,149.43505859375,141727683,2024,8,26,0,0,0,35,0,0,0),
    @(45537,149.6199951171875,160.6600036621094,146.0099945068359,152.8899993896484,152.7010192871094,232644519,2024,9,2,0,0,0,36,0,0,1),
    @(45544,152.6999969482422,156.3500061035156,146.25,150.7200012207031,150.5337066650391,128318915,2024,9,9,0,0,0,37,0,0,0),
    @(45551,151.1999969482422,163.3999938964844,148.5500030517578,161.8099975585938,161.6100006103516,228244249,2024,9,16,0,0,0,38,0,0,0),
    @(45558,166,171,150.3099975585938,151,150.8133544921875,160461262,2024,9,23,0,0,0,39,1,0,0),
    @(45565,150.5,152.7400054931641,137.3399963378906,139.8300018310547,139.8300018310547,88599988,2024,9,30,0,0,0,40,0,0,0),
    @(45572,139.9799957275391,142.8500061035156,128.8099975585938,137.5200042724609,137.5200042724609,107119784,2024,10,7,0,0,0,41,0,0,0),
    @(45579,138.3999938964844,140.7700042724609,126.9800033569336,130.0700073242188,130.0700073242188,62049311,2024,10,14,0,0,0,42,0,0,0),
    @(45586,131.1999969482422,132.8399963378906,111.5100021362305,113.2799987792969,113.2799987792969,122051487,2024,10,21,0,0,0,43,2,0,0),
    @(45593,114.9000015258789,122.8000030517578,112.3499984741211,121.8899993896484,121.8899993896484,62760420,2024,10,28,0,0,0,44,0,0,0),
    @(45600,121.5999984741211,131.5,117.3499984741211,126.7900009155273,126.7900009155273,118877378,2024,11,4,0,0,0,45,0,0,2),
    @(45607,126.0500030517578,133.8000030517578,120.0199966430664,126.8899993896484,126.8899993896484,103738430,2024,11,11,0,0,0,46,0,0,0),
    @(45614,128.6999969482422,132.3999938964844,125.1600036621094,126.2799987792969,126.2799987792969,58458803,2024,11,18,0,0,0,47,0,0,0),
    @(45621,129.3000030517578,135.9299926757812,126.3000030517578,129.0899963378906,129.0899963378906,78236046,2024,11,25,0,0,0,48,0,0,0)
)

$startRow = 1144
for ($i = 0; $i -lt $newRows.Count; $i++) {
  $r = $startRow + $i
  $rowVals = $newRows[$i]

  # Column A carries the datetime style (s="2") used throughout column A;
  # copy it from the row above so formatting stays consistent.
  $ws.Cells.Item($r - 1, 1).Copy($ws.Cells.Item($r, 1))

  for ($c = 0; $c -lt $rowVals.Count; $c++) {
    $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
  }
}
